# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Austria Bundesliga")

function Set-Row($Sheet, $Row, $Values) {
    foreach ($col in $Values.Keys) {
        $Sheet.Range($col + $Row).Value = $Values[$col]
    }
}

# Row 131 <- gets the old content of row 132 (B..AC, column A/C/D/E/H unchanged)
Set-Row $ws 131 @{
    "B"  = 6847108
    "F"  = "Rapid Vienna"
    "G"  = "Austria Lustenau"
    "I"  = 1
    "J"  = "D"
    "K"  = 1.285
    "L"  = 5.5
    "M"  = 8
    "N"  = 1.166
    "O"  = 7
    "P"  = 17
    "Q"  = -2.25
    "R"  = 2.05
    "S"  = 1.8
    "T"  = 3.25
    "U"  = 1.875
    "V"  = 1.975
    "W"  = -1
    "X"  = 6
    "AA" = 0.8
    "AC" = 0.9750000000000001
}

# Row 132 <- gets the old content of row 131
Set-Row $ws 132 @{
    "B"  = 6846476
    "F"  = "FC Salzburg"
    "G"  = "Austria Klagenfurt"
    "I"  = 0
    "J"  = "H"
    "K"  = 1.363
    "L"  = 5
    "M"  = 6.5
    "N"  = 1.285
    "O"  = 5.5
    "P"  = 12
    "Q"  = -1.5
    "R"  = 1.825
    "S"  = 2.025
    "T"  = 2.75
    "U"  = 1.9
    "V"  = 1.95
    "W"  = 0.2849999999999999
    "X"  = -1
    "AA" = 1.025
    "AC" = 0.95
}

# Row 135 <- gets the old content of row 139
Set-Row $ws 135 @{
    "B"  = 6847111
    "F"  = "Wolfsberger AC"
    "G"  = "SCR Altach"
    "H"  = 1
    "I"  = 1
    "J"  = "D"
    "K"  = 2.05
    "L"  = 3.5
    "M"  = 3.5
    "N"  = 2.625
    "O"  = 3.2
    "P"  = 2.75
    "Q"  = 0
    "R"  = 1.875
    "S"  = 1.975
    "T"  = 2
    "U"  = 1.75
    "V"  = 2.05
    "W"  = -1
    "X"  = 2.2
    "Z"  = 0
    "AA" = -0
    "AB" = 0
    "AC" = -0
}

# Row 136 <- gets the old content of row 135
Set-Row $ws 136 @{
    "B"  = 6851939
    "F"  = "Austria Lustenau"
    "G"  = "FC Blau Weiss Linz"
    "H"  = 2
    "I"  = 0
    "J"  = "H"
    "K"  = 4
    "L"  = 3.75
    "M"  = 1.833
    "N"  = 3.75
    "O"  = 3.5
    "P"  = 2
    "Q"  = 0.5
    "R"  = 1.825
    "S"  = 2.025
    "T"  = 2.25
    "U"  = 2
    "V"  = 1.85
    "W"  = 2.75
    "X"  = -1
    "Z"  = 0.825
    "AA" = -1
    "AC" = 0.425
}

# Row 137 <- gets the old content of row 136
Set-Row $ws 137 @{
    "B"  = 6847114
    "F"  = "Hartberg"
    "G"  = "SK Sturm Graz"
    "K"  = 3.4
    "L"  = 3.6
    "M"  = 2
    "N"  = 3.3
    "O"  = 3.3
    "P"  = 2.2
    "Q"  = 0.25
    "R"  = 1.925
    "S"  = 1.925
    "U"  = 1.9
    "V"  = 1.95
    "X"  = 2.3
    "Z"  = 0.4625
    "AC" = 0.475
}

# Row 138 <- gets the old content of row 137
Set-Row $ws 138 @{
    "B"  = 6847113
    "F"  = "Austria Klagenfurt"
    "G"  = "Rapid Vienna"
    "H"  = 1
    "I"  = 1
    "J"  = "D"
    "K"  = 3.6
    "L"  = 3.8
    "M"  = 1.909
    "N"  = 3.6
    "O"  = 2.75
    "P"  = 2.375
    "Q"  = 0.25
    "R"  = 1.85
    "S"  = 2
    "T"  = 2.25
    "U"  = 1.875
    "V"  = 1.975
    "W"  = -1
    "X"  = 1.75
    "Z"  = 0.425
    "AA" = -0.5
    "AB" = -0.5
    "AC" = 0.4875
}

# Row 139 <- gets the old content of row 138
Set-Row $ws 139 @{
    "B"  = 6847112
    "F"  = "FK Austria Vienna"
    "G"  = "WSG Swarovski Tirol"
    "H"  = 2
    "I"  = 0
    "J"  = "H"
    "K"  = 1.5
    "L"  = 4.333
    "M"  = 6
    "N"  = 1.5
    "O"  = 4.5
    "P"  = 6.5
    "Q"  = -1.25
    "R"  = 2.05
    "S"  = 1.8
    "T"  = 2.75
    "U"  = 1.9
    "V"  = 1.95
    "W"  = 0.5
    "X"  = -1
    "Z"  = 1.05
    "AA" = -1
    "AB" = -1
    "AC" = 0.95
}

# Row 152 - targeted updates
Set-Row $ws 152 @{
    "O" = 3.5
    "P" = 4.75
}

# Row 153 - targeted updates
Set-Row $ws 153 @{
    "P" = 10
    "R" = 2.025
    "S" = 1.825
    "U" = 1.9
    "V" = 1.95
}

# Row 154 - targeted updates
Set-Row $ws 154 @{
    "N" = 2.75
    "O" = 3.2
    "P" = 2.7
    "U" = 2.05
    "V" = 1.8
}

# Row 155 - targeted updates
Set-Row $ws 155 @{
    "N" = 2.6
    "O" = 3.25
    "P" = 2.8
}

# Row 156 - targeted updates
Set-Row $ws 156 @{
    "P" = 3.75
    "Q" = -0.5
    "R" = 2.05
    "S" = 1.8
    "U" = 1.85
    "V" = 2
}

# Row 157 - targeted updates
Set-Row $ws 157 @{
    "N" = 1.615
    "R" = 2.1
    "S" = 1.775
    "U" = 2.05
    "V" = 1.8
}
